$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Ridership")

# Update "Riders" (C) and "Average" (D) columns with new Madigan bike hours data
$ws.Range("C2").Value = 91
$ws.Range("D2").Value = 229.3

$ws.Range("C3").Value = 158
$ws.Range("D3").Value = 208.6

$ws.Range("C4").Value = 209
$ws.Range("D4").Value = 193.5

$ws.Range("C5").Value = 247
$ws.Range("D5").Value = 217.78

$ws.Range("C6").Value = 266
$ws.Range("D6").Value = 231.3

$ws.Range("C7").Value = 136
$ws.Range("D7").Value = 118.27

$ws.Range("C8").Value = 86
$ws.Range("D8").Value = 103.33
